$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the existing "_GoBack" bookmark that currently sits at the end of
#    the "Provide any references relevant to the report" paragraph.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2. Fix the misspelling "Receipes" -> "Recipes" in the second paragraph, and
#    drop the spell-check proofErr wrapper that surrounded the old word.
#    (The trailing "isthereanyreceipes.com" URL text must stay untouched.)
# ---------------------------------------------------------------------------
$pTitle = $d.Paragraphs.Item(2)
$titleXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' +
    '<w:p w:rsidR="00EE69F9" w:rsidRDefault="00EE69F9" w:rsidP="00EE69F9">' +
    '<w:r><w:t xml:space="preserve">Is There Any </w:t></w:r>' +
    '<w:r><w:t>Recipes</w:t></w:r>' +
    '<w:r><w:t>? (isthereanyreceipes.com)</w:t></w:r>' +
    '</w:p>' +
    '</w:body></w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'
$pTitle.Range.InsertXML($titleXml)

# ---------------------------------------------------------------------------
# 3. Split the "...You can also construct a prototype..." run in the
#    Navigation Structure flow bullet, and drop a fresh "_GoBack" bookmark
#    right after "You can also" (before the following space).
# ---------------------------------------------------------------------------
$pFlow = $d.Paragraphs.Item(10)
$flowXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' +
    '<w:p w:rsidR="00EE69F9" w:rsidRDefault="00EE69F9" w:rsidP="00EE69F9">' +
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
    '<w:r><w:t>Briefly indicate a typical flow of your application in terms of user experience. You can use any way of representing the flow. You can also</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
    '<w:bookmarkEnd w:id="0"/>' +
    '<w:r><w:t xml:space="preserve"> construct a prototype using one of the prototyping tools to illustrate this. </w:t></w:r>' +
    '</w:p>' +
    '</w:body></w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'
$pFlow.Range.InsertXML($flowXml)
